$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 586 (the "「驚きの価格で」" entry) — everything below shifts up by one.
$ws.Rows.Item(586).Delete()
